$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before row 835, pushing existing rows 835-898 down to 837-900.
$ws.Range("A835:A836").EntireRow.Insert()

# Row 835: new Granny Smith / Primera record
$ws.Cells.Item(835, 1).Value = 5
$ws.Cells.Item(835, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(835, 3).Value = "Maule"
$ws.Cells.Item(835, 4).Value = Get-Date -Year 2022 -Month 5 -Day 25 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(835, 5).Value = 7
$ws.Cells.Item(835, 6).Value = "Fruta"
$ws.Cells.Item(835, 7).Value = 100104
$ws.Cells.Item(835, 8).Value = "Frutos de pepita"
$ws.Cells.Item(835, 9).Value = 100104002
$ws.Cells.Item(835, 10).Value = "Manzana"
$ws.Cells.Item(835, 11).Value = "Granny Smith"
$ws.Cells.Item(835, 12).Value = "Primera"
$ws.Cells.Item(835, 13).Value = 250
$ws.Cells.Item(835, 14).Value = 7000
$ws.Cells.Item(835, 15).Value = 7000
$ws.Cells.Item(835, 16).Value = 7000
$ws.Cells.Item(835, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(835, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(835, 19).Value = 467
$ws.Cells.Item(835, 20).Value = 15

# Row 836: new Pink Lady / Primera record
$ws.Cells.Item(836, 1).Value = 5
$ws.Cells.Item(836, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(836, 3).Value = "Maule"
$ws.Cells.Item(836, 4).Value = Get-Date -Year 2022 -Month 5 -Day 25 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(836, 5).Value = 7
$ws.Cells.Item(836, 6).Value = "Fruta"
$ws.Cells.Item(836, 7).Value = 100104
$ws.Cells.Item(836, 8).Value = "Frutos de pepita"
$ws.Cells.Item(836, 9).Value = 100104002
$ws.Cells.Item(836, 10).Value = "Manzana"
$ws.Cells.Item(836, 11).Value = "Pink Lady"
$ws.Cells.Item(836, 12).Value = "Primera"
$ws.Cells.Item(836, 13).Value = 230
$ws.Cells.Item(836, 14).Value = 7000
$ws.Cells.Item(836, 15).Value = 7000
$ws.Cells.Item(836, 16).Value = 7000
$ws.Cells.Item(836, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(836, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(836, 19).Value = 467
$ws.Cells.Item(836, 20).Value = 15

Write-Output "done"
